# Reverse the order of the comma-separated "Recorded By" entries in column G.
# Only rows whose value actually contains multiple comma-separated parts
# are affected; single-value cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ","
        if ($parts.Count -gt 1) {
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }
            $reversed = @()
            for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
                $reversed += $trimmed[$i]
            }
            $newVal = [string]::Join(", ", $reversed)
            $cell.Value2 = $newVal
        }
    }
}
